# Update the time_taken timestamps on the "data" sheet (column F, rows 2-9)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("F2").Value = "2021-10-05 14:35:41.197334"
$ws.Range("F3").Value = "2021-10-05 14:35:41.197342"
$ws.Range("F4").Value = "2021-10-05 14:35:41.197346"
$ws.Range("F5").Value = "2021-10-05 14:35:41.197348"
$ws.Range("F6").Value = "2021-10-05 14:35:41.197352"
$ws.Range("F7").Value = "2021-10-05 14:35:41.197354"
$ws.Range("F8").Value = "2021-10-05 14:35:41.197357"
$ws.Range("F9").Value = "2021-10-05 14:35:41.197359"

# Add a new "metadata" worksheet, positioned right after "data"
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Match the "data" sheet's outline settings (summary rows below, summary
# columns to the right) and page margins.
$meta.Outline.SummaryRow = 1
$meta.Outline.SummaryColumn = 1

$meta.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$meta.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$meta.PageSetup.TopMargin = $excel.InchesToPoints(1)
$meta.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$meta.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$meta.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# Copy the header/body cell formatting from the "data" sheet so the new
# sheet matches its look (bold, bordered, centered header style; bordered
# style on the first data column).
$ws.Range("B1:G1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Schwannomatosis"
$meta.Range("C2").Value = 3142

# data_version ("0.15") must be stored as text, not a number, so force a
# text number-format before assigning, then drop back to the default
# (unstyled) cell style to match the rest of the row.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.15"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2020-04-24T07:05:46.786724Z"
$meta.Range("F2").Value = "2021-10-05 14:35:41.193500"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3142/?format=json"

# Leave the "data" sheet as the active sheet/tab.
$ws.Activate()
